$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '30.344.07'
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  -2.78%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.942.59'
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -2.80%  '
$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '251.10'
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  -2.47%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.7223'
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  -7.82%  '
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3371'
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  -4.78%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '28.71'
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -2.60%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.07430'
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +5.44%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.8180'
$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  -6.14%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.08143'
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -0.75%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '1.939.16'
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -3.03%  '
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  -1.29%  '
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -5.58%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '14.90'
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  -3.88%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '30.357.53'
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  -2.77%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.000008302'
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  +4.38%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '254.71'
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  -7.23%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '5.895'
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  -1.05%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '2.193.90'
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  -2.83%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.9997'
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '0.9989'
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '6.966'
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  -2.64%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '9.885'
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  -2.28%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '160.60'
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  -2.40%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '2.462'
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  +3.29%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '19.44'
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -2.94%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '0.1321'
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  -11.31%  '
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -2.32%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.349'
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -0.41%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '4.497'
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -2.71%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '4.250'
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  -4.43%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.05266'
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +0.81%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.277'
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  +3.71%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.7575'
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -2.84%  '
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -2.22%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.01993'
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -1.03%  '
$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -2.62%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '82.04'
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  +2.73%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '6.564'
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  -2.90%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.4597'
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -2.87%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.032'
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -5.83%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.8491'
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '103.17'
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -3.15%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '9.915'
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -0.55%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '7.468'
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -3.70%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '37.13'
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  +0.71%  '
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -3.11%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.511'
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -0.37%  '
